# Weekly data refresh: insert two new observation rows (week of 2023-03-23)
# immediately above the former row 1158, which pushes all of the existing
# rows 1158..1221 down to 1160..1223 (same values/format, just relocated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 1158.
$ws.Rows.Item(1158).Insert()
$ws.Rows.Item(1158).Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112021
$categoria = "Ají"
$clasif    = "Hortaliza"

# New row 1158: Ají Americana (o), Primera
$ws.Cells.Item(1158, 1).Value  = $mercadoId
$ws.Cells.Item(1158, 2).Value  = $mercado
$ws.Cells.Item(1158, 3).Value  = $region
$ws.Cells.Item(1158, 4).Value  = 45008
$ws.Cells.Item(1158, 5).Value  = $codreg
$ws.Cells.Item(1158, 6).Value  = $catId
$ws.Cells.Item(1158, 7).Value  = $categoria
$ws.Cells.Item(1158, 8).Value  = "Americana (o)"
$ws.Cells.Item(1158, 9).Value  = "Primera"
$ws.Cells.Item(1158, 10).Value = 250
$ws.Cells.Item(1158, 11).Value = 18000
$ws.Cells.Item(1158, 12).Value = 20000
$ws.Cells.Item(1158, 13).Value = 18800
$ws.Cells.Item(1158, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(1158, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1158, 16).Value = 752
$ws.Cells.Item(1158, 17).Value = 25
$ws.Cells.Item(1158, 18).Value = $clasif

# New row 1159: Ají Americana (o), Segunda
$ws.Cells.Item(1159, 1).Value  = $mercadoId
$ws.Cells.Item(1159, 2).Value  = $mercado
$ws.Cells.Item(1159, 3).Value  = $region
$ws.Cells.Item(1159, 4).Value  = 45008
$ws.Cells.Item(1159, 5).Value  = $codreg
$ws.Cells.Item(1159, 6).Value  = $catId
$ws.Cells.Item(1159, 7).Value  = $categoria
$ws.Cells.Item(1159, 8).Value  = "Americana (o)"
$ws.Cells.Item(1159, 9).Value  = "Segunda"
$ws.Cells.Item(1159, 10).Value = 130
$ws.Cells.Item(1159, 11).Value = 16000
$ws.Cells.Item(1159, 12).Value = 16000
$ws.Cells.Item(1159, 13).Value = 16000
$ws.Cells.Item(1159, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(1159, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1159, 16).Value = 640
$ws.Cells.Item(1159, 17).Value = 25
$ws.Cells.Item(1159, 18).Value = $clasif
